# Update the public EPEX Spot prices workbook:
#  - "Prix Spot": append a new day column T (03-jul) with its 24 hourly prices
#  - "Gaz": append the 2025-07-01 row
#  - "CO2": append the 2025-07-01 row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": append column T for 03-jul
# ---------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the previous day's header cell (S1) into T1 so the new header
# inherits the exact same (bold/bordered/centered) style, then overwrite
# its text with the new date label.
$wsSpot.Range("S1").Copy($wsSpot.Range("T1"))
$wsSpot.Range("T1").Value = "03-jul"

$spotValues = @(
    85,
    81.62,
    79,
    78.93000000000001,
    80.88,
    86.45,
    94.83,
    101.77,
    104.9,
    90.45,
    75.06,
    57.35,
    42.16,
    34.82,
    36.84,
    56.67,
    64.70999999999999,
    81,
    96.34,
    109,
    118.8,
    105.87,
    106.8,
    86.56999999999999
)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 20).Value = $spotValues[$i]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append row 17 for 2025-07-01
# ---------------------------------------------------------------------
# The date is stored as plain text (e.g. "2025-07-01"), not a real Excel
# date, in this workbook. Assigning the string straight to .Value would
# get auto-converted to a date serial, so we build it as a text formula
# first and then flatten it to a static value, which keeps it a plain
# string cell with no special number formatting (matching the rest of
# the column).
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A17").Formula = '="2025-07-01"'
$wsGaz.Range("A17").Copy()
$wsGaz.Range("A17").PasteSpecial(-4163)  # xlPasteValues
$wsGaz.Range("B17").Value = 32.95

# ---------------------------------------------------------------------
# Sheet "CO2": append row 17 for 2025-07-01
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A17").Formula = '="2025-07-01"'
$wsCo2.Range("A17").Copy()
$wsCo2.Range("A17").PasteSpecial(-4163)  # xlPasteValues
$wsCo2.Range("B17").Value = 69.36
